# Auto-generated script: apply 2025-06-28 violent crime daily update
$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet index 1) - 11 cell(s)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 3201  # L2: 3185 -> 3201
$ws.Cells.Item(3, 12).Value = 3290  # L3: 3267 -> 3290
$ws.Cells.Item(4, 9).Value = 1837  # I4: 1836 -> 1837
$ws.Cells.Item(4, 11).Value = 1766  # K4: 1765 -> 1766
$ws.Cells.Item(4, 12).Value = 832  # L4: 825 -> 832
$ws.Cells.Item(5, 12).Value = 186  # L5: 184 -> 186
$ws.Cells.Item(6, 11).Value = 9119  # K6: 9121 -> 9119
$ws.Cells.Item(6, 12).Value = 2908  # L6: 2888 -> 2908
$ws.Cells.Item(7, 9).Value = 26305  # I7: 26304 -> 26305
$ws.Cells.Item(7, 11).Value = 27555  # K7: 27556 -> 27555
$ws.Cells.Item(7, 12).Value = 10417  # L7: 10349 -> 10417

# By Neighborhood (sheet index 2) - 47 cell(s)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(6, 12).Value = 81  # L6: 80 -> 81
$ws.Cells.Item(7, 12).Value = 349  # L7: 347 -> 349
$ws.Cells.Item(8, 12).Value = 663  # L8: 660 -> 663
$ws.Cells.Item(9, 12).Value = 66  # L9: 65 -> 66
$ws.Cells.Item(11, 12).Value = 172  # L11: 171 -> 172
$ws.Cells.Item(13, 12).Value = 15  # L13: 14 -> 15
$ws.Cells.Item(19, 12).Value = 295  # L19: 291 -> 295
$ws.Cells.Item(20, 12).Value = 262  # L20: 261 -> 262
$ws.Cells.Item(24, 12).Value = 24  # L24: 23 -> 24
$ws.Cells.Item(25, 12).Value = 54  # L25: 53 -> 54
$ws.Cells.Item(29, 12).Value = 566  # L29: 562 -> 566
$ws.Cells.Item(31, 12).Value = 100  # L31: 98 -> 100
$ws.Cells.Item(33, 12).Value = 478  # L33: 476 -> 478
$ws.Cells.Item(37, 12).Value = 378  # L37: 373 -> 378
$ws.Cells.Item(40, 12).Value = 26  # L40: 25 -> 26
$ws.Cells.Item(42, 12).Value = 333  # L42: 332 -> 333
$ws.Cells.Item(48, 12).Value = 141  # L48: 140 -> 141
$ws.Cells.Item(49, 12).Value = 57  # L49: 56 -> 57
$ws.Cells.Item(51, 12).Value = 127  # L51: 126 -> 127
$ws.Cells.Item(52, 12).Value = 207  # L52: 206 -> 207
$ws.Cells.Item(53, 12).Value = 119  # L53: 118 -> 119
$ws.Cells.Item(54, 12).Value = 218  # L54: 217 -> 218
$ws.Cells.Item(55, 12).Value = 99  # L55: 98 -> 99
$ws.Cells.Item(63, 9).Value = 261  # I63: 260 -> 261
$ws.Cells.Item(63, 11).Value = 159  # K63: 158 -> 159
$ws.Cells.Item(63, 12).Value = 35  # L63: 33 -> 35
$ws.Cells.Item(64, 12).Value = 69  # L64: 67 -> 69
$ws.Cells.Item(65, 12).Value = 201  # L65: 199 -> 201
$ws.Cells.Item(66, 11).Value = 81  # K66: 82 -> 81
$ws.Cells.Item(67, 12).Value = 382  # L67: 379 -> 382
$ws.Cells.Item(68, 12).Value = 34  # L68: 32 -> 34
$ws.Cells.Item(73, 12).Value = 89  # L73: 90 -> 89
$ws.Cells.Item(76, 12).Value = 142  # L76: 140 -> 142
$ws.Cells.Item(78, 12).Value = 127  # L78: 125 -> 127
$ws.Cells.Item(79, 12).Value = 267  # L79: 266 -> 267
$ws.Cells.Item(82, 11).Value = 32  # K82: 33 -> 32
$ws.Cells.Item(83, 12).Value = 245  # L83: 243 -> 245
$ws.Cells.Item(85, 12).Value = 529  # L85: 526 -> 529
$ws.Cells.Item(86, 12).Value = 78  # L86: 77 -> 78
$ws.Cells.Item(88, 12).Value = 123  # L88: 121 -> 123
$ws.Cells.Item(90, 12).Value = 98  # L90: 96 -> 98
$ws.Cells.Item(91, 12).Value = 152  # L91: 150 -> 152
$ws.Cells.Item(95, 12).Value = 143  # L95: 142 -> 143
$ws.Cells.Item(99, 12).Value = 178  # L99: 176 -> 178
$ws.Cells.Item(101, 9).Value = 26305  # I101: 26304 -> 26305
$ws.Cells.Item(101, 11).Value = 27555  # K101: 27556 -> 27555
$ws.Cells.Item(101, 12).Value = 10417  # L101: 10349 -> 10417

# Auburn Gresham (sheet index 5) - 2 cell(s)
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 12).Value = 110  # L2: 108 -> 110
$ws.Cells.Item(7, 12).Value = 349  # L7: 347 -> 349

# Belmont Cragin (sheet index 6) - 2 cell(s)
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(3, 12).Value = 54  # L3: 53 -> 54
$ws.Cells.Item(7, 12).Value = 172  # L7: 171 -> 172

# South Shore (sheet index 8) - 4 cell(s)
$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 12).Value = 218  # L3: 217 -> 218
$ws.Cells.Item(5, 12).Value = 10  # L5: 9 -> 10
$ws.Cells.Item(6, 12).Value = 108  # L6: 107 -> 108
$ws.Cells.Item(7, 12).Value = 529  # L7: 526 -> 529

# Little Village (sheet index 9) - 2 cell(s)
$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(3, 12).Value = 63  # L3: 62 -> 63
$ws.Cells.Item(7, 12).Value = 207  # L7: 206 -> 207

# Logan Square (sheet index 11) - 2 cell(s)
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(3, 12).Value = 30  # L3: 29 -> 30
$ws.Cells.Item(7, 12).Value = 119  # L7: 118 -> 119

# Austin (sheet index 12) - 3 cell(s)
$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(4, 12).Value = 45  # L4: 44 -> 45
$ws.Cells.Item(6, 12).Value = 187  # L6: 185 -> 187
$ws.Cells.Item(7, 12).Value = 663  # L7: 660 -> 663

# South Chicago (sheet index 13) - 2 cell(s)
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(6, 12).Value = 56  # L6: 54 -> 56
$ws.Cells.Item(7, 12).Value = 245  # L7: 243 -> 245

# Garfield Park (sheet index 14) - 3 cell(s)
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 12).Value = 135  # L2: 134 -> 135
$ws.Cells.Item(6, 12).Value = 165  # L6: 164 -> 165
$ws.Cells.Item(7, 12).Value = 478  # L7: 476 -> 478

# West Pullman (sheet index 15) - 2 cell(s)
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(6, 12).Value = 29  # L6: 28 -> 29
$ws.Cells.Item(7, 12).Value = 143  # L7: 142 -> 143

# Grand Crossing (sheet index 16) - 3 cell(s)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 12).Value = 115  # L2: 112 -> 115
$ws.Cells.Item(3, 12).Value = 112  # L3: 110 -> 112
$ws.Cells.Item(7, 12).Value = 378  # L7: 373 -> 378

# New City (sheet index 17) - 2 cell(s)
$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(6, 12).Value = 57  # L6: 55 -> 57
$ws.Cells.Item(7, 12).Value = 201  # L7: 199 -> 201

# Woodlawn (sheet index 18) - 2 cell(s)
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 12).Value = 74  # L3: 72 -> 74
$ws.Cells.Item(7, 12).Value = 178  # L7: 176 -> 178

# Gage Park (sheet index 20) - 2 cell(s)
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 12).Value = 31  # L2: 29 -> 31
$ws.Cells.Item(7, 12).Value = 100  # L7: 98 -> 100

# North Lawndale (sheet index 21) - 3 cell(s)
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 12).Value = 112  # L2: 111 -> 112
$ws.Cells.Item(6, 12).Value = 86  # L6: 84 -> 86
$ws.Cells.Item(7, 12).Value = 382  # L7: 379 -> 382

# Lincoln Park (sheet index 23) - 2 cell(s)
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(6, 12).Value = 24  # L6: 23 -> 24
$ws.Cells.Item(7, 12).Value = 57  # L7: 56 -> 57

# Loop (sheet index 24) - 2 cell(s)
$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(4, 12).Value = 17  # L4: 16 -> 17
$ws.Cells.Item(7, 12).Value = 218  # L7: 217 -> 218

# Englewood (sheet index 25) - 4 cell(s)
$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 12).Value = 170  # L2: 169 -> 170
$ws.Cells.Item(3, 12).Value = 213  # L3: 211 -> 213
$ws.Cells.Item(6, 12).Value = 148  # L6: 147 -> 148
$ws.Cells.Item(7, 12).Value = 566  # L7: 562 -> 566

# Lake View (sheet index 26) - 2 cell(s)
$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(6, 12).Value = 59  # L6: 58 -> 59
$ws.Cells.Item(7, 12).Value = 141  # L7: 140 -> 141

# Chatham (sheet index 27) - 5 cell(s)
$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 12).Value = 102  # L2: 101 -> 102
$ws.Cells.Item(3, 12).Value = 87  # L3: 86 -> 87
$ws.Cells.Item(4, 12).Value = 12  # L4: 11 -> 12
$ws.Cells.Item(6, 12).Value = 92  # L6: 91 -> 92
$ws.Cells.Item(7, 12).Value = 295  # L7: 291 -> 295

# River North (sheet index 29) - 2 cell(s)
$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(4, 12).Value = 20  # L4: 18 -> 20
$ws.Cells.Item(7, 12).Value = 142  # L7: 140 -> 142

# Ashburn (sheet index 30) - 2 cell(s)
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(6, 12).Value = 18  # L6: 17 -> 18
$ws.Cells.Item(7, 12).Value = 81  # L7: 80 -> 81

# Humboldt Park (sheet index 32) - 2 cell(s)
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 12).Value = 100  # L2: 99 -> 100
$ws.Cells.Item(7, 12).Value = 333  # L7: 332 -> 333

# Boystown (sheet index 33) - 2 cell(s)
$ws = $wb.Worksheets.Item('Boystown')
$ws.Cells.Item(5, 12).Value = 7  # L5: 6 -> 7
$ws.Cells.Item(6, 12).Value = 15  # L6: 14 -> 15

# Rogers Park (sheet index 35) - 3 cell(s)
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 12).Value = 37  # L2: 36 -> 37
$ws.Cells.Item(3, 12).Value = 37  # L3: 36 -> 37
$ws.Cells.Item(7, 12).Value = 127  # L7: 125 -> 127

# Lower West Side (sheet index 36) - 2 cell(s)
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 12).Value = 34  # L3: 33 -> 34
$ws.Cells.Item(7, 12).Value = 99  # L7: 98 -> 99

# Dunning (sheet index 37) - 2 cell(s)
$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(2, 12).Value = 13  # L2: 12 -> 13
$ws.Cells.Item(7, 12).Value = 24  # L7: 23 -> 24

# Washington Park (sheet index 40) - 3 cell(s)
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(3, 12).Value = 62  # L3: 61 -> 62
$ws.Cells.Item(5, 12).Value = 5  # L5: 4 -> 5
$ws.Cells.Item(7, 12).Value = 152  # L7: 150 -> 152

# Roseland (sheet index 42) - 2 cell(s)
$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 12).Value = 94  # L3: 93 -> 94
$ws.Cells.Item(7, 12).Value = 267  # L7: 266 -> 267

# Near South Side (sheet index 43) - 2 cell(s)
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 12).Value = 18  # L3: 16 -> 18
$ws.Cells.Item(7, 12).Value = 69  # L7: 67 -> 69

# Chicago Lawn (sheet index 44) - 2 cell(s)
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 12).Value = 79  # L3: 78 -> 79
$ws.Cells.Item(7, 12).Value = 262  # L7: 261 -> 262

# East Side (sheet index 52) - 2 cell(s)
$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(3, 12).Value = 28  # L3: 27 -> 28
$ws.Cells.Item(7, 12).Value = 54  # L7: 53 -> 54

# North Center (sheet index 59) - 2 cell(s)
$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 11).Value = 38  # K6: 39 -> 38
$ws.Cells.Item(7, 11).Value = 81  # K7: 82 -> 81

# Avalon Park (sheet index 61) - 2 cell(s)
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(3, 12).Value = 27  # L3: 26 -> 27
$ws.Cells.Item(7, 12).Value = 66  # L7: 65 -> 66

# Portage Park (sheet index 62) - 2 cell(s)
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(6, 12).Value = 21  # L6: 22 -> 21
$ws.Cells.Item(7, 12).Value = 89  # L7: 90 -> 89

# West Town (sheet index 65) - 2 cell(s)
$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 12).Value = 18  # L3: 17 -> 18
$ws.Cells.Item(6, 12).Value = 50  # L6: 51 -> 50

# United Center (sheet index 68) - 2 cell(s)
$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 12).Value = 38  # L6: 36 -> 38
$ws.Cells.Item(7, 12).Value = 123  # L7: 121 -> 123

# Streeterville (sheet index 72) - 2 cell(s)
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(3, 12).Value = 15  # L3: 14 -> 15
$ws.Cells.Item(7, 12).Value = 78  # L7: 77 -> 78

# Washington Heights (sheet index 74) - 3 cell(s)
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(3, 12).Value = 28  # L3: 27 -> 28
$ws.Cells.Item(5, 12).Value = 2  # L5: 1 -> 2
$ws.Cells.Item(7, 12).Value = 98  # L7: 96 -> 98

# Little Italy, UIC (sheet index 75) - 2 cell(s)
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(4, 12).Value = 18  # L4: 17 -> 18
$ws.Cells.Item(7, 12).Value = 127  # L7: 126 -> 127

# North Park (sheet index 76) - 3 cell(s)
$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(2, 12).Value = 9  # L2: 8 -> 9
$ws.Cells.Item(6, 12).Value = 11  # L6: 10 -> 11
$ws.Cells.Item(7, 12).Value = 34  # L7: 32 -> 34

# Sheffield & DePaul (sheet index 83) - 2 cell(s)
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(6, 11).Value = 17  # K6: 18 -> 17
$ws.Cells.Item(7, 11).Value = 32  # K7: 33 -> 32

# Hegewisch (sheet index 89) - 2 cell(s)
$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(6, 12).Value = 6  # L6: 5 -> 6
$ws.Cells.Item(7, 12).Value = 26  # L7: 25 -> 26
